$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Cells.Item(28, 8).Value = 1059.6316
$wb.Worksheets.Item(1).Cells.Item(28, 9).Value = 1257.091
$wb.Worksheets.Item(1).Cells.Item(28, 10).Value = 788.125
$wb.Worksheets.Item(1).Cells.Item(28, 11).Value = 1257.091
$wb.Worksheets.Item(1).Cells.Item(28, 12).Value = 788.125
$wb.Worksheets.Item(1).Cells.Item(28, 13).Value = -772.0909999999999
$wb.Worksheets.Item(1).Cells.Item(28, 14).Value = -1758.125
$wb.Worksheets.Item(1).Cells.Item(70, 8).Value = 711.53845
$wb.Worksheets.Item(1).Cells.Item(70, 10).Value = 711.53845
$wb.Worksheets.Item(1).Cells.Item(70, 12).Value = 2134.61535
$wb.Worksheets.Item(1).Cells.Item(70, 14).Value = -2674.61535
$wb.Worksheets.Item(1).Cells.Item(73, 8).Value = 711.53845
$wb.Worksheets.Item(1).Cells.Item(73, 10).Value = 711.53845
$wb.Worksheets.Item(1).Cells.Item(73, 12).Value = 2134.61535
$wb.Worksheets.Item(1).Cells.Item(73, 14).Value = -4006.61535
$wb.Worksheets.Item(1).Cells.Item(132, 8).Value = 906537.0600000001
$wb.Worksheets.Item(1).Cells.Item(132, 9).Value = 1059908.6
$wb.Worksheets.Item(1).Cells.Item(132, 10).Value = 1645
$wb.Worksheets.Item(1).Cells.Item(132, 11).Value = 3179725.8
$wb.Worksheets.Item(1).Cells.Item(132, 12).Value = 4935
$wb.Worksheets.Item(1).Cells.Item(132, 13).Value = -3177195.8
$wb.Worksheets.Item(1).Cells.Item(132, 14).Value = -9995
$wb.Worksheets.Item(1).Cells.Item(137, 8).Value = 857.7778
$wb.Worksheets.Item(1).Cells.Item(137, 9).Value = 754.8276
$wb.Worksheets.Item(1).Cells.Item(137, 10).Value = 1284.2858
$wb.Worksheets.Item(1).Cells.Item(137, 11).Value = 2264.4828
$wb.Worksheets.Item(1).Cells.Item(137, 12).Value = 3852.8574
$wb.Worksheets.Item(1).Cells.Item(137, 13).Value = 285.5172000000002
$wb.Worksheets.Item(1).Cells.Item(137, 14).Value = -8952.857400000001
$wb.Worksheets.Item(1).Cells.Item(138, 8).Value = 3791.6836
$wb.Worksheets.Item(1).Cells.Item(138, 9).Value = 1086.5
$wb.Worksheets.Item(1).Cells.Item(138, 10).Value = 5820.5713
$wb.Worksheets.Item(1).Cells.Item(138, 11).Value = 3259.5
$wb.Worksheets.Item(1).Cells.Item(138, 12).Value = 17461.7139
$wb.Worksheets.Item(1).Cells.Item(138, 13).Value = 1880.5
$wb.Worksheets.Item(1).Cells.Item(138, 14).Value = -27741.7139
$wb.Worksheets.Item(2).Cells.Item(31, 8).Value = 22190.334
$wb.Worksheets.Item(2).Cells.Item(31, 9).Value = 6735.5
$wb.Worksheets.Item(2).Cells.Item(31, 10).Value = 53100
$wb.Worksheets.Item(2).Cells.Item(31, 11).Value = 6735.5
$wb.Worksheets.Item(2).Cells.Item(31, 12).Value = 53100
$wb.Worksheets.Item(2).Cells.Item(31, 13).Value = -6441.5
$wb.Worksheets.Item(2).Cells.Item(31, 14).Value = -53688
$wb.Worksheets.Item(2).Cells.Item(97, 8).Value = 598.3333
$wb.Worksheets.Item(2).Cells.Item(97, 9).Value = 318
$wb.Worksheets.Item(2).Cells.Item(97, 10).Value = 2000
$wb.Worksheets.Item(2).Cells.Item(97, 11).Value = 318
$wb.Worksheets.Item(2).Cells.Item(97, 12).Value = 2000
$wb.Worksheets.Item(2).Cells.Item(97, 13).Value = 178
$wb.Worksheets.Item(2).Cells.Item(97, 14).Value = -2992
$wb.Worksheets.Item(2).Cells.Item(132, 8).Value = 1476.566
$wb.Worksheets.Item(2).Cells.Item(132, 9).Value = 1219.8536
$wb.Worksheets.Item(2).Cells.Item(132, 10).Value = 2353.6667
$wb.Worksheets.Item(2).Cells.Item(132, 11).Value = 3659.5608
$wb.Worksheets.Item(2).Cells.Item(132, 12).Value = 7061.000100000001
$wb.Worksheets.Item(2).Cells.Item(132, 13).Value = -1129.5608
$wb.Worksheets.Item(2).Cells.Item(132, 14).Value = -12121.0001
$wb.Worksheets.Item(3).Cells.Item(134, 8).Value = 4727.5713
$wb.Worksheets.Item(3).Cells.Item(134, 9).Value = 5264
$wb.Worksheets.Item(3).Cells.Item(134, 10).Value = 2260
$wb.Worksheets.Item(3).Cells.Item(134, 11).Value = 15792
$wb.Worksheets.Item(3).Cells.Item(134, 12).Value = 6780
$wb.Worksheets.Item(3).Cells.Item(134, 13).Value = -13257
$wb.Worksheets.Item(3).Cells.Item(134, 14).Value = -11850
$wb.Worksheets.Item(4).Cells.Item(31, 8).Value = 2926166.8
$wb.Worksheets.Item(4).Cells.Item(31, 9).Value = 3269660
$wb.Worksheets.Item(4).Cells.Item(31, 10).Value = 6474.75
$wb.Worksheets.Item(4).Cells.Item(31, 11).Value = 3269660
$wb.Worksheets.Item(4).Cells.Item(31, 12).Value = 6474.75
$wb.Worksheets.Item(4).Cells.Item(31, 13).Value = -3269365
$wb.Worksheets.Item(4).Cells.Item(31, 14).Value = -7064.75
$wb.Worksheets.Item(4).Cells.Item(34, 8).Value = 2926166.8
$wb.Worksheets.Item(4).Cells.Item(34, 9).Value = 3269660
$wb.Worksheets.Item(4).Cells.Item(34, 10).Value = 6474.75
$wb.Worksheets.Item(4).Cells.Item(34, 11).Value = 3269660
$wb.Worksheets.Item(4).Cells.Item(34, 12).Value = 6474.75
$wb.Worksheets.Item(4).Cells.Item(34, 13).Value = -3269458
$wb.Worksheets.Item(4).Cells.Item(34, 14).Value = -6878.75
$wb.Worksheets.Item(4).Cells.Item(58, 8).Value = 719153.5
$wb.Worksheets.Item(4).Cells.Item(58, 9).Value = 906415.5
$wb.Worksheets.Item(4).Cells.Item(58, 10).Value = 1315.8334
$wb.Worksheets.Item(4).Cells.Item(58, 11).Value = 906415.5
$wb.Worksheets.Item(4).Cells.Item(58, 12).Value = 1315.8334
$wb.Worksheets.Item(4).Cells.Item(58, 13).Value = -906212.5
$wb.Worksheets.Item(4).Cells.Item(58, 14).Value = -1721.8334
$wb.Worksheets.Item(4).Cells.Item(132, 8).Value = 1659.9412
$wb.Worksheets.Item(4).Cells.Item(132, 9).Value = 1298.1034
$wb.Worksheets.Item(4).Cells.Item(132, 10).Value = 3758.6
$wb.Worksheets.Item(4).Cells.Item(132, 11).Value = 3894.3102
$wb.Worksheets.Item(4).Cells.Item(132, 12).Value = 11275.8
$wb.Worksheets.Item(4).Cells.Item(132, 13).Value = -1364.3102
$wb.Worksheets.Item(4).Cells.Item(132, 14).Value = -16335.8
$wb.Worksheets.Item(4).Cells.Item(134, 8).Value = 1558.4546
$wb.Worksheets.Item(4).Cells.Item(134, 9).Value = 1338.6154
$wb.Worksheets.Item(4).Cells.Item(134, 10).Value = 1876
$wb.Worksheets.Item(4).Cells.Item(134, 11).Value = 4015.8462
$wb.Worksheets.Item(4).Cells.Item(134, 12).Value = 5628
$wb.Worksheets.Item(4).Cells.Item(134, 13).Value = -1480.8462
$wb.Worksheets.Item(4).Cells.Item(134, 14).Value = -10698
$wb.Worksheets.Item(4).Cells.Item(136, 8).Value = 719153.5
$wb.Worksheets.Item(4).Cells.Item(136, 9).Value = 906415.5
$wb.Worksheets.Item(4).Cells.Item(136, 10).Value = 1315.8334
$wb.Worksheets.Item(4).Cells.Item(136, 11).Value = 2719246.5
$wb.Worksheets.Item(4).Cells.Item(136, 12).Value = 3947.5002
$wb.Worksheets.Item(4).Cells.Item(136, 13).Value = -2716696.5
$wb.Worksheets.Item(4).Cells.Item(136, 14).Value = -9047.5002
$wb.Worksheets.Item(5).Cells.Item(5, 8).Value = 528.3333
$wb.Worksheets.Item(5).Cells.Item(5, 9).Value = 416.3889
$wb.Worksheets.Item(5).Cells.Item(5, 10).Value = 1200
$wb.Worksheets.Item(5).Cells.Item(5, 11).Value = 1249.1667
$wb.Worksheets.Item(5).Cells.Item(5, 12).Value = 3600
$wb.Worksheets.Item(5).Cells.Item(5, 13).Value = -1137.1667
$wb.Worksheets.Item(5).Cells.Item(5, 14).Value = -3824
$wb.Worksheets.Item(5).Cells.Item(132, 8).Value = 2144.1428
$wb.Worksheets.Item(5).Cells.Item(132, 9).Value = 1502
$wb.Worksheets.Item(5).Cells.Item(132, 10).Value = 2401
$wb.Worksheets.Item(5).Cells.Item(132, 11).Value = 13518
$wb.Worksheets.Item(5).Cells.Item(132, 12).Value = 21609
$wb.Worksheets.Item(5).Cells.Item(132, 13).Value = -10988
$wb.Worksheets.Item(5).Cells.Item(132, 14).Value = -26669
$wb.Worksheets.Item(5).Cells.Item(135, 8).Value = 528.3333
$wb.Worksheets.Item(5).Cells.Item(135, 9).Value = 416.3889
$wb.Worksheets.Item(5).Cells.Item(135, 10).Value = 1200
$wb.Worksheets.Item(5).Cells.Item(135, 11).Value = 3747.5001
$wb.Worksheets.Item(5).Cells.Item(135, 12).Value = 10800
$wb.Worksheets.Item(5).Cells.Item(135, 13).Value = -1212.5001
$wb.Worksheets.Item(5).Cells.Item(135, 14).Value = -15870
$wb.Worksheets.Item(5).Cells.Item(138, 8).Value = 17859988
$wb.Worksheets.Item(5).Cells.Item(138, 9).Value = 19233580
$wb.Worksheets.Item(5).Cells.Item(138, 11).Value = 57700740
$wb.Worksheets.Item(5).Cells.Item(138, 13).Value = -57695600
$wb.Worksheets.Item(7).Cells.Item(7, 8).Value = 2272.6155
$wb.Worksheets.Item(7).Cells.Item(7, 9).Value = 2271.6428
$wb.Worksheets.Item(7).Cells.Item(7, 10).Value = 2273.75
$wb.Worksheets.Item(7).Cells.Item(7, 11).Value = 2271.6428
$wb.Worksheets.Item(7).Cells.Item(7, 12).Value = 2273.75
$wb.Worksheets.Item(7).Cells.Item(7, 13).Value = -2159.6428
$wb.Worksheets.Item(7).Cells.Item(7, 14).Value = -2497.75
$wb.Worksheets.Item(7).Cells.Item(22, 8).Value = 536.6429000000001
$wb.Worksheets.Item(7).Cells.Item(22, 9).Value = 501.18182
$wb.Worksheets.Item(7).Cells.Item(22, 10).Value = 666.6667
$wb.Worksheets.Item(7).Cells.Item(22, 11).Value = 501.18182
$wb.Worksheets.Item(7).Cells.Item(22, 12).Value = 666.6667
$wb.Worksheets.Item(7).Cells.Item(22, 13).Value = -206.18182
$wb.Worksheets.Item(7).Cells.Item(22, 14).Value = -1256.6667
$wb.Worksheets.Item(7).Cells.Item(27, 8).Value = 536.6429000000001
$wb.Worksheets.Item(7).Cells.Item(27, 9).Value = 501.18182
$wb.Worksheets.Item(7).Cells.Item(27, 10).Value = 666.6667
$wb.Worksheets.Item(7).Cells.Item(27, 11).Value = 501.18182
$wb.Worksheets.Item(7).Cells.Item(27, 12).Value = 666.6667
$wb.Worksheets.Item(7).Cells.Item(27, 13).Value = -394.18182
$wb.Worksheets.Item(7).Cells.Item(27, 14).Value = -880.6667
$wb.Worksheets.Item(7).Cells.Item(68, 8).Value = 14222.5
$wb.Worksheets.Item(7).Cells.Item(68, 9).Value = 26075
$wb.Worksheets.Item(7).Cells.Item(68, 10).Value = 2370
$wb.Worksheets.Item(7).Cells.Item(68, 11).Value = 26075
$wb.Worksheets.Item(7).Cells.Item(68, 12).Value = 2370
$wb.Worksheets.Item(7).Cells.Item(68, 13).Value = -25326
$wb.Worksheets.Item(7).Cells.Item(68, 14).Value = -3868
$wb.Worksheets.Item(7).Cells.Item(71, 8).Value = 14222.5
$wb.Worksheets.Item(7).Cells.Item(71, 9).Value = 26075
$wb.Worksheets.Item(7).Cells.Item(71, 10).Value = 2370
$wb.Worksheets.Item(7).Cells.Item(71, 11).Value = 130375
$wb.Worksheets.Item(7).Cells.Item(71, 12).Value = 11850
$wb.Worksheets.Item(7).Cells.Item(71, 13).Value = -126631
$wb.Worksheets.Item(7).Cells.Item(71, 14).Value = -19338
$wb.Worksheets.Item(7).Cells.Item(122, 8).Value = 2930
$wb.Worksheets.Item(7).Cells.Item(122, 9).Value = 2540
$wb.Worksheets.Item(7).Cells.Item(122, 10).Value = 3320
$wb.Worksheets.Item(7).Cells.Item(122, 11).Value = 7620
$wb.Worksheets.Item(7).Cells.Item(122, 12).Value = 9960
$wb.Worksheets.Item(7).Cells.Item(122, 13).Value = -5170
$wb.Worksheets.Item(7).Cells.Item(122, 14).Value = -14860
$wb.Worksheets.Item(7).Cells.Item(126, 8).Value = 2272.6155
$wb.Worksheets.Item(7).Cells.Item(126, 9).Value = 2271.6428
$wb.Worksheets.Item(7).Cells.Item(126, 10).Value = 2273.75
$wb.Worksheets.Item(7).Cells.Item(126, 11).Value = 6814.928400000001
$wb.Worksheets.Item(7).Cells.Item(126, 12).Value = 6821.25
$wb.Worksheets.Item(7).Cells.Item(126, 13).Value = -4344.928400000001
$wb.Worksheets.Item(7).Cells.Item(126, 14).Value = -11761.25
$wb.Worksheets.Item(7).Cells.Item(127, 8).Value = 35000
$wb.Worksheets.Item(7).Cells.Item(127, 10).Value = 35000
$wb.Worksheets.Item(7).Cells.Item(127, 12).Value = 35000
$wb.Worksheets.Item(7).Cells.Item(127, 14).Value = -44920
$wb.Worksheets.Item(7).Cells.Item(132, 8).Value = 3510.7273
$wb.Worksheets.Item(7).Cells.Item(132, 9).Value = 1795.5
$wb.Worksheets.Item(7).Cells.Item(132, 10).Value = 4490.857
$wb.Worksheets.Item(7).Cells.Item(132, 11).Value = 5386.5
$wb.Worksheets.Item(7).Cells.Item(132, 12).Value = 13472.571
$wb.Worksheets.Item(7).Cells.Item(132, 13).Value = -2856.5
$wb.Worksheets.Item(7).Cells.Item(132, 14).Value = -18532.571
$wb.Worksheets.Item(7).Cells.Item(136, 8).Value = 2054.5356
$wb.Worksheets.Item(7).Cells.Item(136, 9).Value = 1101.6666
$wb.Worksheets.Item(7).Cells.Item(136, 10).Value = 2769.1875
$wb.Worksheets.Item(7).Cells.Item(136, 11).Value = 3304.9998
$wb.Worksheets.Item(7).Cells.Item(136, 12).Value = 8307.5625
$wb.Worksheets.Item(7).Cells.Item(136, 13).Value = -754.9998000000001
$wb.Worksheets.Item(7).Cells.Item(136, 14).Value = -13407.5625
$wb.Worksheets.Item(8).Cells.Item(113, 8).Value = 363.26923
$wb.Worksheets.Item(8).Cells.Item(113, 9).Value = 269.6
$wb.Worksheets.Item(8).Cells.Item(113, 10).Value = 491
$wb.Worksheets.Item(8).Cells.Item(113, 11).Value = 808.8000000000001
$wb.Worksheets.Item(8).Cells.Item(113, 12).Value = 1473
$wb.Worksheets.Item(8).Cells.Item(113, 13).Value = 1361.2
$wb.Worksheets.Item(8).Cells.Item(113, 14).Value = -5813
$wb.Worksheets.Item(8).Cells.Item(126, 8).Value = 1025835.44
$wb.Worksheets.Item(8).Cells.Item(126, 9).Value = 1367121.9
$wb.Worksheets.Item(8).Cells.Item(126, 10).Value = 1976.25
$wb.Worksheets.Item(8).Cells.Item(126, 11).Value = 4101365.7
$wb.Worksheets.Item(8).Cells.Item(126, 12).Value = 5928.75
$wb.Worksheets.Item(8).Cells.Item(126, 13).Value = -4098895.7
$wb.Worksheets.Item(8).Cells.Item(126, 14).Value = -10868.75
$wb.Worksheets.Item(8).Cells.Item(131, 8).Value = 45375.125
$wb.Worksheets.Item(8).Cells.Item(131, 9).Value = 0
$wb.Worksheets.Item(8).Cells.Item(131, 10).Value = 45375.125
$wb.Worksheets.Item(8).Cells.Item(131, 11).Value = 0
$wb.Worksheets.Item(8).Cells.Item(131, 12).Value = 45375.125
$wb.Worksheets.Item(8).Cells.Item(131, 13).ClearContents()
$wb.Worksheets.Item(8).Cells.Item(131, 14).Value = -55455.125
$wb.Worksheets.Item(8).Cells.Item(132, 8).Value = 1666.675
$wb.Worksheets.Item(8).Cells.Item(132, 9).Value = 1332.4584
$wb.Worksheets.Item(8).Cells.Item(132, 10).Value = 2168
$wb.Worksheets.Item(8).Cells.Item(132, 11).Value = 3997.3752
$wb.Worksheets.Item(8).Cells.Item(132, 12).Value = 6504
$wb.Worksheets.Item(8).Cells.Item(132, 13).Value = -1467.3752
$wb.Worksheets.Item(8).Cells.Item(132, 14).Value = -11564
$wb.Worksheets.Item(8).Cells.Item(136, 8).Value = 1701564
$wb.Worksheets.Item(8).Cells.Item(136, 9).Value = 2778494
$wb.Worksheets.Item(8).Cells.Item(136, 10).Value = 1148.421
$wb.Worksheets.Item(8).Cells.Item(136, 11).Value = 8335482
$wb.Worksheets.Item(8).Cells.Item(136, 12).Value = 3445.263
$wb.Worksheets.Item(8).Cells.Item(136, 13).Value = -8332932
$wb.Worksheets.Item(8).Cells.Item(136, 14).Value = -8545.262999999999
